# Fix curly brace replacement: the "Report date" template cell used a
# format specifier Excel can't evaluate ("{Date: MMM dd, yyyy}" - the
# leading space after the colon breaks the placeholder's format-string
# parsing). Replace it with a valid, unambiguous custom date format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Report date: {Date:dd.MM.yyyy}"
